$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1143057577823544
$ws.Range("D2").Value = 0.91003221000996

$ws.Range("C3").Value = -0.09292503697217756
$ws.Range("D3").Value = 0.9268041652979311

$ws.Range("C4").Value = 0.8767942290025227
$ws.Range("D4").Value = 0.390072490932093

$ws.Range("C5").Value = 1.9173645297332
$ws.Range("D5").Value = 0.06826775212219971

$ws.Range("C6").Value = 0.008832211092691418
$ws.Range("D6").Value = 0.9930326098626041

$ws.Range("C7").Value = 0.8752370211093564
$ws.Range("D7").Value = 0.3909006404183586

$ws.Range("C8").Value = 1.907661244122388
$ws.Range("D8").Value = 0.06957327732500462
$ws.Range("G8").Value = "No"

$ws.Range("C9").Value = 0.7019714046853849
$ws.Range("D9").Value = 0.4900575996616339

$ws.Range("C10").Value = 1.341279824576258
$ws.Range("D10").Value = 0.193513207742801

$ws.Range("C11").Value = 0.7894232184578512
$ws.Range("D11").Value = 0.4382905013504033
